$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny float precision drift on the existing A5 timestamp
$ws.Range("A5").Value = 45806.39330972223

# Append new row with the latest price entry
$ws.Range("A6").NumberFormat = $ws.Range("A5").NumberFormat
$ws.Range("A6").Value = 45806.40671058597
$ws.Range("B6").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C6").Value = "1Kg"
$ws.Range("D6").Value = "12,88€"
